$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Swap match data (columns F:V only - match/odds info) between row pairs
# that got reordered in the source feed.
# ---------------------------------------------------------------------------

function Swap-Rows($ws, $rowA, $rowB) {
    $rangeA = $ws.Range("F" + $rowA + ":V" + $rowA)
    $rangeB = $ws.Range("F" + $rowB + ":V" + $rowB)
    $valA = $rangeA.Value2
    $valB = $rangeB.Value2
    $rangeA.Value2 = $valB
    $rangeB.Value2 = $valA
}

# Rows 24 <-> 25
Swap-Rows $ws 24 25

# Rows 41 <-> 42
Swap-Rows $ws 41 42

# Rows 52 -> 53 -> 54 -> 52 (3-way rotation):
#   new52 = old54 ; new53 = old52 ; new54 = old53
$range52 = $ws.Range("F52:V52")
$range53 = $ws.Range("F53:V53")
$range54 = $ws.Range("F54:V54")
$val52 = $range52.Value2
$val53 = $range53.Value2
$val54 = $range54.Value2
$range52.Value2 = $val54
$range53.Value2 = $val52
$range54.Value2 = $val53

# Rows 61 <-> 62
Swap-Rows $ws 61 62

# Rows 77 <-> 78
Swap-Rows $ws 77 78

# ---------------------------------------------------------------------------
# Append three new match rows (86, 87, 88) at the end of the sheet,
# cloning the formatting of the last existing row (85).
# ---------------------------------------------------------------------------

$ws.Range("A85:V85").Copy($ws.Range("A86:V88"))

$ws.Range("A86").Value2 = 85
$ws.Range("B86").Value2 = "montenegro"
$ws.Range("C86").Value2 = "prva-crnogorska-liga"
$ws.Range("D86").Value2 = "2023-2024"
$ws.Range("E86").Value2 = 45262.54166666666
$ws.Range("F86").Value2 = "Decic"
$ws.Range("G86").Value2 = 0
$ws.Range("H86").Value2 = "Arsenal Tivat"
$ws.Range("I86").Value2 = 1
$ws.Range("J86").Value2 = 1.47
$ws.Range("K86").Value2 = "01/12/2023 01:13"
$ws.Range("L86").Value2 = 1.41
$ws.Range("M86").Value2 = "02/12/2023 12:56"
$ws.Range("N86").Value2 = 3.78
$ws.Range("O86").Value2 = "01/12/2023 01:13"
$ws.Range("P86").Value2 = 4.22
$ws.Range("Q86").Value2 = "02/12/2023 12:56"
$ws.Range("R86").Value2 = 6.2
$ws.Range("S86").Value2 = "01/12/2023 01:13"
$ws.Range("T86").Value2 = 8.050000000000001
$ws.Range("U86").Value2 = "02/12/2023 12:56"
$ws.Range("V86").Value2 = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/decic-arsenal-tivat/jPojAg5T/"

$ws.Range("A87").Value2 = 86
$ws.Range("B87").Value2 = "montenegro"
$ws.Range("C87").Value2 = "prva-crnogorska-liga"
$ws.Range("D87").Value2 = "2023-2024"
$ws.Range("E87").Value2 = 45262.54166666666
$ws.Range("F87").Value2 = "Rudar"
$ws.Range("G87").Value2 = 1
$ws.Range("H87").Value2 = "Sutjeska"
$ws.Range("I87").Value2 = 2
$ws.Range("J87").Value2 = 5.42
$ws.Range("K87").Value2 = "01/12/2023 01:13"
$ws.Range("L87").Value2 = 4.66
$ws.Range("M87").Value2 = "02/12/2023 12:40"
$ws.Range("N87").Value2 = 3.57
$ws.Range("O87").Value2 = "01/12/2023 01:13"
$ws.Range("P87").Value2 = 3.26
$ws.Range("Q87").Value2 = "02/12/2023 12:54"
$ws.Range("R87").Value2 = 1.56
$ws.Range("S87").Value2 = "01/12/2023 01:13"
$ws.Range("T87").Value2 = 1.82
$ws.Range("U87").Value2 = "02/12/2023 12:54"
$ws.Range("V87").Value2 = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/rudar-sutjeska/CYnnBZkN/"

$ws.Range("A88").Value2 = 87
$ws.Range("B88").Value2 = "montenegro"
$ws.Range("C88").Value2 = "prva-crnogorska-liga"
$ws.Range("D88").Value2 = "2023-2024"
$ws.Range("E88").Value2 = 45262.625
$ws.Range("F88").Value2 = "Mladost DG"
$ws.Range("G88").Value2 = 2
$ws.Range("H88").Value2 = "Jedinstvo"
$ws.Range("I88").Value2 = 0
$ws.Range("J88").Value2 = 2.14
$ws.Range("K88").Value2 = "01/12/2023 03:12"
$ws.Range("L88").Value2 = 2.19
$ws.Range("M88").Value2 = "02/12/2023 14:37"
$ws.Range("N88").Value2 = 3
$ws.Range("O88").Value2 = "01/12/2023 03:12"
$ws.Range("P88").Value2 = 2.9
$ws.Range("Q88").Value2 = "02/12/2023 14:37"
$ws.Range("R88").Value2 = 3.27
$ws.Range("S88").Value2 = "01/12/2023 03:12"
$ws.Range("T88").Value2 = 3.7
$ws.Range("U88").Value2 = "02/12/2023 14:37"
$ws.Range("V88").Value2 = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/mladost-dg-jedinstvo/84YH5DDp/"
